$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.248195171356201
$ws.Range("B1").Value = 2.645231246948242
$ws.Range("C1").Value = 8.423727989196777
$ws.Range("D1").Value = 2.110546827316284
$ws.Range("E1").Value = 1.138992428779602
